$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay as text, matching original inlineStr cells
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.637.42"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "2.626.35"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "524.31"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "143.28"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("E9").Value = "  -7.47%  "
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "3.086.66"
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("D14").Value = "58.554.86"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").Value = "21.10"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "2.538.20"
$ws.Range("E17").Value = "  -7.62%  "
$ws.Range("D18").Value = "339.11"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "65.36"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "7.18"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").Value = "6.52"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "150.14"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("D36").Value = "0.900"
$ws.Range("D37").Value = "0.855"
$ws.Range("E37").Value = "  -5.80%  "
$ws.Range("D38").Value = "36.41"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("E39").Value = "  -6.84%  "
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").Value = "0.0971"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "270.76"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").Value = "0.0533"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "19.13"
$ws.Range("E47").Value = "  -5.80%  "
$ws.Range("D48").Value = "2.037.57"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("E50").Value = "  -8.09%  "
$ws.Range("D51").Value = "18.39"
$ws.Range("E51").Value = "  -5.66%  "
